# Generate Report for Handoff
# Adds two new tracked files (9b93defb-... and 9c3fb80a-...) to the
# localization-status workbook: one new row per file on the "Overview"
# sheet, and one new row per file on each language sheet ("zh-cn", "de-de").

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Helper data for the two new handed-off files
# ---------------------------------------------------------------------
$uuidA = "9b93defb-8a5e-4e70-8f2d-5227b63eee93"
$uuidB = "9c3fb80a-df32-4a39-8ba1-8174f4fca193"

$hashA = "49aa8aa2ca5bca26ec18ed37110b1ba3a74df61c"
$hashB = "4ef33e3323be7380751e7462bc673b318cd7b194"

$mdA = "$uuidA.md"
$mdB = "$uuidB.md"

$statusReady = "Ready for handoff"
$dateA = "2016-03-21 20:37:49"
$dateZhA = "2016-03-21 20:37:45"
$emptyHandback = "0001-01-01 00:00:00"
$include = "Include"
$ext = ".md"

$xlfZhA = "$uuidA.$hashA.zh-cn.xlf"
$xlfZhB = "$uuidB.$hashB.zh-cn.xlf"
$xlfDeA = "$uuidA.$hashA.de-de.xlf"
$xlfDeB = "$uuidB.$hashB.de-de.xlf"

$mdUrlA = "https://github.com/OpenLocalizationTest/oltest/blob/0000000000000000000000000000000000000000/e2e/$mdA"
$mdUrlB = "https://github.com/OpenLocalizationTest/oltest/blob/0000000000000000000000000000000000000000/e2e/$mdB"

$xlfUrlZhA = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0000000000000000000000000000000000000000/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$xlfZhA"
$xlfUrlZhB = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0000000000000000000000000000000000000000/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$xlfZhB"
$xlfUrlDeA = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0000000000000000000000000000000000000000/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$xlfDeA"
$xlfUrlDeB = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0000000000000000000000000000000000000000/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$xlfDeB"

# ---------------------------------------------------------------------
# Sheet 1 : "Overview" -> rows 4 & 5 (A:D), hyperlinks on column A
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Cells.Item(4, 2).Value = $statusReady
$wsOverview.Cells.Item(4, 3).Value = $statusReady
$wsOverview.Cells.Item(4, 4).Value = $dateA
$wsOverview.Cells.Item(4, 4).NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsOverview.Cells.Item(5, 2).Value = $statusReady
$wsOverview.Cells.Item(5, 3).Value = $statusReady
$wsOverview.Cells.Item(5, 4).Value = $dateA
$wsOverview.Cells.Item(5, 4).NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsOverview.Cells.Item(4, 1).Value = $mdA
$wsOverview.Hyperlinks.Add($wsOverview.Range("A4"), $mdUrlA, "", "", $mdA) | Out-Null
$wsOverview.Cells.Item(4, 1).Font.Underline = 2
$wsOverview.Cells.Item(4, 1).Font.Color = 6591981

$wsOverview.Cells.Item(5, 1).Value = $mdB
$wsOverview.Hyperlinks.Add($wsOverview.Range("A5"), $mdUrlB, "", "", $mdB) | Out-Null
$wsOverview.Cells.Item(5, 1).Font.Underline = 2
$wsOverview.Cells.Item(5, 1).Font.Color = 6591981

# ---------------------------------------------------------------------
# Sheet 2 : "zh-cn" -> rows 4 & 5
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Cells.Item(4, 2).Value = $ext
$wsZh.Cells.Item(4, 3).Value = $statusReady
$wsZh.Cells.Item(4, 5).Value = $dateZhA
$wsZh.Cells.Item(4, 5).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Cells.Item(4, 8).Value = $emptyHandback
$wsZh.Cells.Item(4, 8).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Cells.Item(4, 10).Value = $include

$wsZh.Cells.Item(5, 2).Value = $ext
$wsZh.Cells.Item(5, 3).Value = $statusReady
$wsZh.Cells.Item(5, 5).Value = $dateZhA
$wsZh.Cells.Item(5, 5).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Cells.Item(5, 8).Value = $emptyHandback
$wsZh.Cells.Item(5, 8).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Cells.Item(5, 10).Value = $include

$wsZh.Cells.Item(4, 1).Value = $mdA
$wsZh.Hyperlinks.Add($wsZh.Range("A4"), $mdUrlA, "", "", $mdA) | Out-Null
$wsZh.Cells.Item(4, 1).Font.Underline = 2
$wsZh.Cells.Item(4, 1).Font.Color = 6591981

$wsZh.Cells.Item(4, 4).Value = $xlfZhA
$wsZh.Hyperlinks.Add($wsZh.Range("D4"), $xlfUrlZhA, "", "", $xlfZhA) | Out-Null
$wsZh.Cells.Item(4, 4).Font.Underline = 2
$wsZh.Cells.Item(4, 4).Font.Color = 6591981

$wsZh.Cells.Item(5, 1).Value = $mdB
$wsZh.Hyperlinks.Add($wsZh.Range("A5"), $mdUrlB, "", "", $mdB) | Out-Null
$wsZh.Cells.Item(5, 1).Font.Underline = 2
$wsZh.Cells.Item(5, 1).Font.Color = 6591981

$wsZh.Cells.Item(5, 4).Value = $xlfZhB
$wsZh.Hyperlinks.Add($wsZh.Range("D5"), $xlfUrlZhB, "", "", $xlfZhB) | Out-Null
$wsZh.Cells.Item(5, 4).Font.Underline = 2
$wsZh.Cells.Item(5, 4).Font.Color = 6591981

# ---------------------------------------------------------------------
# Sheet 3 : "de-de" -> rows 4 & 5
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Cells.Item(4, 2).Value = $ext
$wsDe.Cells.Item(4, 3).Value = $statusReady
$wsDe.Cells.Item(4, 5).Value = $dateA
$wsDe.Cells.Item(4, 5).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Cells.Item(4, 8).Value = $emptyHandback
$wsDe.Cells.Item(4, 8).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Cells.Item(4, 10).Value = $include

$wsDe.Cells.Item(5, 2).Value = $ext
$wsDe.Cells.Item(5, 3).Value = $statusReady
$wsDe.Cells.Item(5, 5).Value = $dateA
$wsDe.Cells.Item(5, 5).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Cells.Item(5, 8).Value = $emptyHandback
$wsDe.Cells.Item(5, 8).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Cells.Item(5, 10).Value = $include

$wsDe.Cells.Item(4, 1).Value = $mdA
$wsDe.Hyperlinks.Add($wsDe.Range("A4"), $mdUrlA, "", "", $mdA) | Out-Null
$wsDe.Cells.Item(4, 1).Font.Underline = 2
$wsDe.Cells.Item(4, 1).Font.Color = 6591981

$wsDe.Cells.Item(4, 4).Value = $xlfDeA
$wsDe.Hyperlinks.Add($wsDe.Range("D4"), $xlfUrlDeA, "", "", $xlfDeA) | Out-Null
$wsDe.Cells.Item(4, 4).Font.Underline = 2
$wsDe.Cells.Item(4, 4).Font.Color = 6591981

$wsDe.Cells.Item(5, 1).Value = $mdB
$wsDe.Hyperlinks.Add($wsDe.Range("A5"), $mdUrlB, "", "", $mdB) | Out-Null
$wsDe.Cells.Item(5, 1).Font.Underline = 2
$wsDe.Cells.Item(5, 1).Font.Color = 6591981

$wsDe.Cells.Item(5, 4).Value = $xlfDeB
$wsDe.Hyperlinks.Add($wsDe.Range("D5"), $xlfUrlDeB, "", "", $xlfDeB) | Out-Null
$wsDe.Cells.Item(5, 4).Font.Underline = 2
$wsDe.Cells.Item(5, 4).Font.Color = 6591981

Write-Host "Generate Report for Handoff: added rows 4-5 to Overview, zh-cn, de-de sheets."
